$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Rules")

# Row 10 ("Integer min" / R40 rule): "From" bound (C10) changes from 18 to 1
$ws.Range("C10").Value = 1
